$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.606.07"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "2.559.48"
$ws.Range("E3").Value = "  -3.24%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'512.18"
$ws.Range("E5").Value = "  -2.83%  "
$ws.Range("D6").Value = "'139.61"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.559"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").Value = "2.567.57"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "'6.40"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").Value = "'0.0988"
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("D12").Value = "'0.329"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "3.009.19"
$ws.Range("E14").Value = "  -3.17%  "
$ws.Range("D15").Value = "57.528.16"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "'19.93"
$ws.Range("E16").Value = "  -5.27%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.555.95"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000131"
$ws.Range("E18").Value = "  -4.40%  "
$ws.Range("D19").Value = "'331.09"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").Value = "'4.23"
$ws.Range("E20").Value = "  -5.28%  "
$ws.Range("D21").Value = "'9.98"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").Value = "'6.28"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'64.60"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "'0.165"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.693.22"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").Value = "'0.394"
$ws.Range("E28").Value = "  -6.56%  "
$ws.Range("D29").Value = "'6.90"
$ws.Range("E29").Value = "  -4.90%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "0.0₃0712"
$ws.Range("E31").Value = "  -11.24%  "
$ws.Range("D32").Value = "'6.01"
$ws.Range("E32").Value = "  -8.21%  "
$ws.Range("D33").Value = "'1.55"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'149.02"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'18.38"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("D36").Value = "'3.88"
$ws.Range("E36").Value = "  -8.16%  "
$ws.Range("D37").Value = "'1.10"
$ws.Range("E37").Value = "  -8.23%  "
$ws.Range("D38").Value = "'0.828"
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("D39").Value = "'35.56"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").Value = "'0.820"
$ws.Range("E40").Value = "  -10.72%  "
$ws.Range("D41").Value = "'1.43"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'0.996"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.42"
$ws.Range("E43").Value = "  -6.48%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "'10.68"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.590"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'264.39"
$ws.Range("E46").Value = "  -3.71%  "
$ws.Range("D47").Value = "'0.0938"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0512"
$ws.Range("E48").Value = "  -5.15%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'18.30"
$ws.Range("E49").Value = "  -5.81%  "
$ws.Range("D50").Value = "1.942.65"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0217"
$ws.Range("E51").Value = "  -6.29%  "

# Reset style on quote-prefixed numeric-looking text cells so no extra
# quotePrefix style attribute lingers on the cell (keeps cells styleless,
# matching the original unstyled inlineStr cells) while preserving the
# literal text value set above.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
